# ---------------------------------------------------------------------------
# Update gh-pages output (合肥-漫展信息.xlsx) to match the refreshed scrape.
#  - refresh "想去人数" (F column) counters on sheets "展览" and "全部类型"
#  - swap in a brand-new event ("合肥·第八届环形宇宙动漫游戏嘉年华Plus") in
#    place of the existing "银魂主题派对only2.0" row, then re-insert the
#    "银魂" event as a new row right after it (with its own refreshed
#    counter), pushing the following rows down by one
# ---------------------------------------------------------------------------

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("展览")
$ws4 = $wb.Worksheets.Item("全部类型")

# --- Update F column (想去人数) counts in sheet '展览' (Exhibition) ---
$ws1.Range("F2").Value = 135
$ws1.Range("F3").Value = 234
$ws1.Range("F4").Value = 16
$ws1.Range("F5").Value = 6713
$ws1.Range("F6").Value = 85
$ws1.Range("F7").Value = 432
$ws1.Range("F8").Value = 139
$ws1.Range("F9").Value = 6232
$ws1.Range("F13").Value = 12
$ws1.Range("F14").Value = 96
$ws1.Range("F16").Value = 122
$ws1.Range("F17").Value = 19
$ws1.Range("F18").Value = 365
$ws1.Range("F19").Value = 44
$ws1.Range("F20").Value = 8
$ws1.Range("F21").Value = 4546
$ws1.Range("F22").Value = 58
$ws1.Range("F23").Value = 38

# --- Update F column (想去人数) counts in sheet '全部类型' (All types) ---
$ws4.Range("F2").Value = 135
$ws4.Range("F3").Value = 234
$ws4.Range("F4").Value = 16
$ws4.Range("F5").Value = 6713
$ws4.Range("F6").Value = 85
$ws4.Range("F7").Value = 432
$ws4.Range("F8").Value = 139
$ws4.Range("F9").Value = 6232
$ws4.Range("F13").Value = 12
$ws4.Range("F14").Value = 96
$ws4.Range("F16").Value = 122
$ws4.Range("F17").Value = 19
$ws4.Range("F18").Value = 365
$ws4.Range("F19").Value = 44
$ws4.Range("F20").Value = 8
$ws4.Range("F21").Value = 4546
$ws4.Range("F23").Value = 58
$ws4.Range("F24").Value = 38

# ---------------------------------------------------------------------------
# Helper: write a text value into a cell without letting Excel's "looks
# like a date" auto-detection turn it into a real date/time value - force
# Text format, assign, then drop back to the sheet's default (General)
# style so the cell ends up looking just like its neighbours.
# ---------------------------------------------------------------------------
function Set-TextValue {
    param($cell, [string]$text)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# ---------------------------------------------------------------------------
# On each sheet:
#   1) the "银魂主题派对only2.0" row is overwritten in place with the new
#      "环形宇宙...Plus" event (its row number / A-column index do not move)
#   2) a new row is inserted right below it, re-creating the old "银魂"
#      event with its own refreshed "想去人数" count - this pushes every
#      following row (starting with "SSS第五人格only") down by one
#   3) the now-shifted "SSS第五人格only" row gets its "想去人数" refreshed
# ---------------------------------------------------------------------------
function Update-YinhunBlock {
    param($ws, [int]$yinhunRow)

    $newEventRow = $yinhunRow
    $insertedYinhunRow = $yinhunRow + 1
    $sssRow = $yinhunRow + 2

    $counterBefore = $ws.Range("A" + $newEventRow).Value()

    # 1) Overwrite the existing "银魂" row in place with the new event.
    $ws.Range("C" + $newEventRow).Value = "合肥·第八届环形宇宙动漫游戏嘉年华Plus"
    $ws.Range("D" + $newEventRow).Value = "南京路与庐州大道交汇处 合肥滨湖国际会展中心"
    $ws.Range("E" + $newEventRow).Value = "2024.08.17 09:30-08.18 17:00"
    $ws.Range("F" + $newEventRow).Value = 10
    $ws.Range("G" + $newEventRow).Value = 69
    $ws.Range("H" + $newEventRow).Value = "https://show.bilibili.com/platform/detail.html?id=88650"
    $ws.Range("I" + $newEventRow).Value = "//i2.hdslb.com/bfs/openplatform/202407/4I7mduRV1720071650216.jpeg"

    # 2) Insert a fresh row right after it and restore the "银魂" event
    #    there, copying column-A's bold/border style down from the row
    #    above so the new row matches the rest of the table.
    $ws.Rows.Item($insertedYinhunRow).Insert()
    $ws.Range("A" + $newEventRow).Copy()
    $ws.Range("A" + $insertedYinhunRow).PasteSpecial(-4122)
    $ws.Range("A" + $insertedYinhunRow).Value = $counterBefore + 1

    Set-TextValue $ws.Range("B" + $insertedYinhunRow) "2024-08-17"
    $ws.Range("C" + $insertedYinhunRow).Value = "合肥·银魂主题派对only2.0"
    $ws.Range("D" + $insertedYinhunRow).Value = "长江东路1137号圣大国际商贸中心2-301室 梦田音乐LiveHouse(合肥店)"
    $ws.Range("E" + $insertedYinhunRow).Value = "2024.08.17 13:00-08.17 18:00"
    $ws.Range("F" + $insertedYinhunRow).Value = 192
    $ws.Range("G" + $insertedYinhunRow).Value = 128
    $ws.Range("H" + $insertedYinhunRow).Value = "https://show.bilibili.com/platform/detail.html?id=87173"
    $ws.Range("I" + $insertedYinhunRow).Value = "//i2.hdslb.com/bfs/openplatform/202406/aSc8SoTl1718078234193.png"

    # 3) The old "SSS第五人格only" row has now shifted down to $sssRow -
    #    just refresh its "想去人数" counter.
    $ws.Range("F" + $sssRow).Value = 66
}

Update-YinhunBlock $ws1 24
Update-YinhunBlock $ws4 25
